# Updates cryptos list values (prices / 1h volume %) per the Fri Jul 28 2023
# GitHub Actions data refresh, including the Filecoin/PancakeSwap row swap.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    # Force text storage so numeric-looking strings (and strings with
    # trailing zeros such as "1.000") are not auto-converted to numbers,
    # then restore the default "Normal" style so no stray number format
    # is left attached to the cell.
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "29.243.25"
Set-TextValue "E2" "  -0.66%  "
Set-TextValue "D3" "1.862.07"
Set-TextValue "E3" "  -0.86%  "
Set-TextValue "D4" "1.000"
Set-TextValue "E4" "  +0.05%  "
Set-TextValue "E5" "  -0.75%  "
Set-TextValue "D6" "240.68"
Set-TextValue "E6" "  +0.20%  "
Set-TextValue "D7" "0.9993"
Set-TextValue "E7" "  -0.04%  "
Set-TextValue "D8" "0.3087"
Set-TextValue "E8" "  -0.59%  "
Set-TextValue "E9" "  -1.52%  "
Set-TextValue "D10" "24.99"
Set-TextValue "E10" "  +0.10%  "
Set-TextValue "D11" "0.08316"
Set-TextValue "E11" "  +0.67%  "
Set-TextValue "D12" "1.870.81"
Set-TextValue "E12" "  +0.32%  "
Set-TextValue "D13" "0.7183"
Set-TextValue "E13" "  -1.34%  "
Set-TextValue "D14" "5.223"
Set-TextValue "E14" "  -1.17%  "
Set-TextValue "D15" "90.88"
Set-TextValue "E15" "  -0.44%  "
Set-TextValue "D16" "29.248.15"
Set-TextValue "E16" "  -0.40%  "
Set-TextValue "D17" "5.992"
Set-TextValue "E17" "  +0.92%  "
Set-TextValue "D18" "243.80"
Set-TextValue "E18" "  -0.64%  "
Set-TextValue "D19" "2.144.61"
Set-TextValue "E19" "  +1.71%  "
Set-TextValue "D20" "0.000007813"
Set-TextValue "E20" "  -1.21%  "
Set-TextValue "D21" "13.17"
Set-TextValue "E21" "  -1.09%  "
Set-TextValue "D22" "0.9990"
Set-TextValue "E22" "  -0.02%  "
Set-TextValue "D23" "7.951"
Set-TextValue "E23" "  +0.02%  "
Set-TextValue "D24" "0.9998"
Set-TextValue "E24" "  +0.03%  "
Set-TextValue "D25" "0.1611"
Set-TextValue "E25" "  +1.48%  "
Set-TextValue "D26" "162.68"
Set-TextValue "E26" "  -0.82%  "
Set-TextValue "D27" "8.918"
Set-TextValue "E27" "  -1.27%  "
Set-TextValue "D28" "18.60"
Set-TextValue "E28" "  +1.40%  "
Set-TextValue "D29" "1.357"
Set-TextValue "E29" "  -0.60%  "
Set-TextValue "B30" "PancakeSwap"
Set-TextValue "C30" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D30" "1.498"
Set-TextValue "E30" "  +0.89%  "
Set-TextValue "B31" "Filecoin"
Set-TextValue "C31" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D31" "4.436"
Set-TextValue "E31" "  +0.89%  "
Set-TextValue "D32" "4.257"
Set-TextValue "E32" "  +2.49%  "
Set-TextValue "D33" "0.05186"
Set-TextValue "E33" "  -1.89%  "
Set-TextValue "D34" "0.8206"
Set-TextValue "E34" "  +13.37%  "
Set-TextValue "D35" "1.934"
Set-TextValue "E35" "  -0.62%  "
Set-TextValue "D36" "1.175"
Set-TextValue "E36" "  -2.23%  "
Set-TextValue "D37" "2.680"
Set-TextValue "E37" "  +0.10%  "
Set-TextValue "E38" "  -0.48%  "
Set-TextValue "D39" "2.693"
Set-TextValue "E39" "  -0.96%  "
Set-TextValue "D40" "1.156.29"
Set-TextValue "E40" "  -6.90%  "
Set-TextValue "D41" "6.215"
Set-TextValue "E41" "  +2.07%  "
Set-TextValue "D42" "0.8985"
Set-TextValue "E42" "  -1.16%  "
Set-TextValue "D43" "72.91"
Set-TextValue "E43" "  -0.43%  "
Set-TextValue "D44" "0.9983"
Set-TextValue "E44" "  -0.16%  "
Set-TextValue "D45" "101.92"
Set-TextValue "E45" "  -1.78%  "
Set-TextValue "D46" "2.042.27"
Set-TextValue "E46" "  +1.50%  "
Set-TextValue "D47" "0.5180"
Set-TextValue "E47" "  -2.86%  "
Set-TextValue "D48" "1.784"
Set-TextValue "E48" "  +1.05%  "
Set-TextValue "D49" "9.380"
Set-TextValue "E49" "  +0.86%  "
Set-TextValue "E50" "  -1.10%  "
Set-TextValue "D51" "7.080"
Set-TextValue "E51" "  -0.09%  "
